$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 13 new rows at the very top; this shifts all existing rows (and their
# formatting / row spans) down by 13, matching rows 1-56 -> 14-69 and blank
# rows 57-61 -> 70-74.
$ws.Rows("1:13").Insert()

# Column A (dates) does not have a column-level style, so copy the date
# number format from the first pre-existing row (now row 14) onto the new
# rows, reusing the existing style instead of creating a new one.
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A1:A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$nbsp = [char]0x00A0

$ws.Range("A1").Value2 = 41746
$ws.Range("B1").Value2 = 'RETIRO ATM BP D/REINA VICTORIA'
$ws.Range("C1").Value2 = 'D'
$ws.Range("D1").Value2 = '0000990490'
$ws.Range("E1").Value2 = 'CENTRO DE ACOPIO NORTE'
$ws.Range("F1").Value2 = [string]::Concat('20.00', $nbsp, $nbsp)
$ws.Range("G1").Value2 = '748.22'

$ws.Range("A2").Value2 = 41744
$ws.Range("B2").Value2 = 'DEP CNB-1500415029001'
$ws.Range("C2").Value2 = 'C'
$ws.Range("D2").Value2 = '0008775881'
$ws.Range("E2").Value2 = 'AG. NORTE'
$ws.Range("F2").Value2 = [string]::Concat('50.00', $nbsp, $nbsp)
$ws.Range("G2").Value2 = '768.22'

$ws.Range("A3").Value2 = 41744
$ws.Range("B3").Value2 = 'CONSUMO DATA FYBECA (PLAZA DE TOROS'
$ws.Range("C3").Value2 = 'D'
$ws.Range("D3").Value2 = '0004392772'
$ws.Range("E3").Value2 = 'INSTITUCIONAL SS.CC.'
$ws.Range("F3").Value2 = [string]::Concat('41.58', $nbsp, $nbsp)
$ws.Range("G3").Value2 = '718.22'

$ws.Range("A4").Value2 = 41744
$ws.Range("B4").Value2 = 'CONSUMO DATA FYBECA TORRES MEDICAS'
$ws.Range("C4").Value2 = 'D'
$ws.Range("D4").Value2 = '0002223892'
$ws.Range("E4").Value2 = 'INSTITUCIONAL SS.CC.'
$ws.Range("F4").Value2 = [string]::Concat('2.31', $nbsp, $nbsp)
$ws.Range("G4").Value2 = '759.80'

$ws.Range("A5").Value2 = 41744
$ws.Range("B5").Value2 = 'CONSUMO DATA FYBECA TORRES MEDICAS'
$ws.Range("C5").Value2 = 'D'
$ws.Range("D5").Value2 = '0002195122'
$ws.Range("E5").Value2 = 'INSTITUCIONAL SS.CC.'
$ws.Range("F5").Value2 = [string]::Concat('38.19', $nbsp, $nbsp)
$ws.Range("G5").Value2 = '762.11'

$ws.Range("A6").Value2 = 41743
$ws.Range("B6").Value2 = '  TRANSFERENCIA INTERNET'
$ws.Range("C6").Value2 = 'D'
$ws.Range("D6").Value2 = '0007431020'
$ws.Range("E6").Value2 = 'AG. NORTE'
$ws.Range("F6").Value2 = [string]::Concat('100.00', $nbsp, $nbsp)
$ws.Range("G6").Value2 = '800.30'

$ws.Range("A7").Value2 = 41743
$ws.Range("B7").Value2 = '  TRANSFERENCIA INTERNET'
$ws.Range("C7").Value2 = 'D'
$ws.Range("D7").Value2 = '0007349236'
$ws.Range("E7").Value2 = 'AG. NORTE'
$ws.Range("F7").Value2 = [string]::Concat('243.19', $nbsp, $nbsp)
$ws.Range("G7").Value2 = '900.30'

$ws.Range("A8").Value2 = 41743
$ws.Range("B8").Value2 = 'CONSUMO DATA AKI MOLINEROS 161'
$ws.Range("C8").Value2 = 'D'
$ws.Range("D8").Value2 = '0002730824'
$ws.Range("E8").Value2 = 'INSTITUCIONAL SS.CC.'
$ws.Range("F8").Value2 = [string]::Concat('85.43', $nbsp, $nbsp)
$ws.Range("G8").Value2 = '1143.49'

$ws.Range("A9").Value2 = 41740
$ws.Range("B9").Value2 = 'CONSUMO DATA FYBECA TORRES MEDICAS'
$ws.Range("C9").Value2 = 'D'
$ws.Range("D9").Value2 = '0004017189'
$ws.Range("E9").Value2 = 'INSTITUCIONAL SS.CC.'
$ws.Range("F9").Value2 = [string]::Concat('3.37', $nbsp, $nbsp)
$ws.Range("G9").Value2 = '1228.92'

$ws.Range("A10").Value2 = 41740
$ws.Range("B10").Value2 = 'RETIRO ATM BP D/H. METROPOLITANO'
$ws.Range("C10").Value2 = 'D'
$ws.Range("D10").Value2 = '0002026113'
$ws.Range("E10").Value2 = 'HOSPITAL METROPOLITANO'
$ws.Range("F10").Value2 = [string]::Concat('200.00', $nbsp, $nbsp)
$ws.Range("G10").Value2 = '1232.29'

$ws.Range("A11").Value2 = 41739
$ws.Range("B11").Value2 = 'CONSUMO DATA FYBECA (PLAZA DE TOROS'
$ws.Range("C11").Value2 = 'D'
$ws.Range("D11").Value2 = '0008586167'
$ws.Range("E11").Value2 = 'INSTITUCIONAL SS.CC.'
$ws.Range("F11").Value2 = [string]::Concat('26.80', $nbsp, $nbsp)
$ws.Range("G11").Value2 = '1432.29'

$ws.Range("A12").Value2 = 41739
$ws.Range("B12").Value2 = 'RETIRO ATM BP D/H. METROPOLITANO'
$ws.Range("C12").Value2 = 'D'
$ws.Range("D12").Value2 = '0008502524'
$ws.Range("E12").Value2 = 'HOSPITAL METROPOLITANO'
$ws.Range("F12").Value2 = [string]::Concat('50.00', $nbsp, $nbsp)
$ws.Range("G12").Value2 = '1459.09'

$ws.Range("A13").Value2 = 41738
$ws.Range("B13").Value2 = '  TRANSFERENCIA INTERNET'
$ws.Range("C13").Value2 = 'D'
$ws.Range("D13").Value2 = '0004566186'
$ws.Range("E13").Value2 = 'AG. NORTE'
$ws.Range("F13").Value2 = [string]::Concat('1000.00', $nbsp, $nbsp)
$ws.Range("G13").Value2 = '1509.09'


function Get-MoFormula($row) {
    return "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A$row,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B$row,""', 'mo_tipo' => '"",C$row,""', 'mo_documento' => '"",D$row,""', 'mo_oficina' => '"",E$row,""', 'mo_monto' => "",TRIM(F$row),"", 'mo_saldo' => "",G$row,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_borrado_logico' => false),"")"
}

for ($r = 1; $r -le 8; $r++) {
    $ws.Range("H$r").Formula = Get-MoFormula $r
}

# Update the selection to mirror the new shared-formula range.
$ws.Range("H1:H8").Select() | Out-Null
